$wb = $excel.ActiveWorkbook

# ===== Sheet: LP1912 =====
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 08:31:01"
$ws1.Range("A3").Value = "Total filas: 116"

# -- changed rows --
$ws1.Cells.Item(26,1).Value = "05:55:02"
$ws1.Cells.Item(26,2).Value = "06:30"
$ws1.Cells.Item(26,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Cells.Item(26,4).Value = 35
$ws1.Cells.Item(26,5).Value = "LP1912"
$ws1.Cells.Item(27,1).Value = "06:25:28"
$ws1.Cells.Item(27,2).Value = "06:30"
$ws1.Cells.Item(27,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(27,4).Value = 5
$ws1.Cells.Item(27,5).Value = "LP1912"
$ws1.Cells.Item(32,1).Value = "06:54:06"
$ws1.Cells.Item(32,2).Value = "06:55"
$ws1.Cells.Item(32,3).Value = "14_ABASTO"
$ws1.Cells.Item(32,4).Value = 1
$ws1.Cells.Item(32,5).Value = "LP1912"
$ws1.Cells.Item(33,1).Value = "06:54:06"
$ws1.Cells.Item(33,2).Value = "06:55"
$ws1.Cells.Item(33,3).Value = "215C_EL PATO"
$ws1.Cells.Item(33,4).Value = 1
$ws1.Cells.Item(33,5).Value = "LP1912"
$ws1.Cells.Item(37,1).Value = "06:25:28"
$ws1.Cells.Item(37,2).Value = "07:05"
$ws1.Cells.Item(37,3).Value = "15_ABASTO"
$ws1.Cells.Item(37,4).Value = 40
$ws1.Cells.Item(37,5).Value = "LP1912"
$ws1.Cells.Item(38,1).Value = "06:54:06"
$ws1.Cells.Item(38,2).Value = "07:05"
$ws1.Cells.Item(38,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(38,4).Value = 11
$ws1.Cells.Item(38,5).Value = "LP1912"
$ws1.Cells.Item(49,1).Value = "07:17:59"
$ws1.Cells.Item(49,2).Value = "07:31"
$ws1.Cells.Item(49,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(49,4).Value = 14
$ws1.Cells.Item(49,5).Value = "LP1912"
$ws1.Cells.Item(50,1).Value = "07:17:59"
$ws1.Cells.Item(50,2).Value = "07:31"
$ws1.Cells.Item(50,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(50,4).Value = 14
$ws1.Cells.Item(50,5).Value = "LP1912"
$ws1.Cells.Item(51,1).Value = "07:17:59"
$ws1.Cells.Item(51,2).Value = "07:31"
$ws1.Cells.Item(51,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(51,4).Value = 14
$ws1.Cells.Item(51,5).Value = "LP1912"
$ws1.Cells.Item(52,1).Value = "05:55:02"
$ws1.Cells.Item(52,2).Value = "07:32"
$ws1.Cells.Item(52,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(52,4).Value = 97
$ws1.Cells.Item(52,5).Value = "LP1912"
$ws1.Cells.Item(53,1).Value = "06:54:06"
$ws1.Cells.Item(53,2).Value = "07:32"
$ws1.Cells.Item(53,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(53,4).Value = 38
$ws1.Cells.Item(53,5).Value = "LP1912"
$ws1.Cells.Item(54,1).Value = "06:54:06"
$ws1.Cells.Item(54,2).Value = "07:32"
$ws1.Cells.Item(54,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(54,4).Value = 38
$ws1.Cells.Item(54,5).Value = "LP1912"
$ws1.Cells.Item(78,1).Value = "07:17:59"
$ws1.Cells.Item(78,2).Value = "08:22"
$ws1.Cells.Item(78,3).Value = "215B_EL PATO"
$ws1.Cells.Item(78,4).Value = 65
$ws1.Cells.Item(78,5).Value = "LP1912"
$ws1.Cells.Item(79,1).Value = "07:17:59"
$ws1.Cells.Item(79,2).Value = "08:22"
$ws1.Cells.Item(79,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(79,4).Value = 65
$ws1.Cells.Item(79,5).Value = "LP1912"
$ws1.Cells.Item(87,1).Value = "08:31:01"
$ws1.Cells.Item(87,2).Value = "08:35"
$ws1.Cells.Item(87,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(87,4).Value = 4
$ws1.Cells.Item(87,5).Value = "LP1912"
$ws1.Cells.Item(91,1).Value = "08:31:01"
$ws1.Cells.Item(91,2).Value = "08:42"
$ws1.Cells.Item(91,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(91,4).Value = 11
$ws1.Cells.Item(91,5).Value = "LP1912"
$ws1.Cells.Item(93,1).Value = "08:31:01"
$ws1.Cells.Item(93,2).Value = "08:44"
$ws1.Cells.Item(93,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(93,4).Value = 13
$ws1.Cells.Item(93,5).Value = "LP1912"
$ws1.Cells.Item(94,1).Value = "08:31:01"
$ws1.Cells.Item(94,2).Value = "08:44"
$ws1.Cells.Item(94,3).Value = "14_ABASTO"
$ws1.Cells.Item(94,4).Value = 13
$ws1.Cells.Item(94,5).Value = "LP1912"
$ws1.Cells.Item(95,1).Value = "08:01:10"
$ws1.Cells.Item(95,2).Value = "08:49"
$ws1.Cells.Item(95,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(95,4).Value = 48
$ws1.Cells.Item(95,5).Value = "LP1912"
$ws1.Cells.Item(96,1).Value = "07:17:59"
$ws1.Cells.Item(96,2).Value = "08:53"
$ws1.Cells.Item(96,3).Value = "17_ROMERO"
$ws1.Cells.Item(96,4).Value = 96
$ws1.Cells.Item(96,5).Value = "LP1912"
$ws1.Cells.Item(97,1).Value = "08:31:01"
$ws1.Cells.Item(97,2).Value = "08:53"
$ws1.Cells.Item(97,3).Value = "10_OLMOS"
$ws1.Cells.Item(97,4).Value = 22
$ws1.Cells.Item(97,5).Value = "LP1912"
$ws1.Cells.Item(98,1).Value = "08:31:01"
$ws1.Cells.Item(98,2).Value = "08:54"
$ws1.Cells.Item(98,3).Value = "17_ROMERO"
$ws1.Cells.Item(98,4).Value = 23
$ws1.Cells.Item(98,5).Value = "LP1912"
$ws1.Cells.Item(99,1).Value = "07:17:59"
$ws1.Cells.Item(99,2).Value = "09:01"
$ws1.Cells.Item(99,3).Value = "215A_EL PATO"
$ws1.Cells.Item(99,4).Value = 104
$ws1.Cells.Item(99,5).Value = "LP1912"
$ws1.Cells.Item(100,1).Value = "08:31:01"
$ws1.Cells.Item(100,2).Value = "09:02"
$ws1.Cells.Item(100,3).Value = "215A_EL PATO"
$ws1.Cells.Item(100,4).Value = 31
$ws1.Cells.Item(100,5).Value = "LP1912"
$ws1.Cells.Item(101,1).Value = "08:31:01"
$ws1.Cells.Item(101,2).Value = "09:04"
$ws1.Cells.Item(101,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(101,4).Value = 33
$ws1.Cells.Item(101,5).Value = "LP1912"
$ws1.Cells.Item(102,1).Value = "08:31:01"
$ws1.Cells.Item(102,2).Value = "09:05"
$ws1.Cells.Item(102,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(102,4).Value = 34
$ws1.Cells.Item(102,5).Value = "LP1912"
$ws1.Cells.Item(103,1).Value = "07:17:59"
$ws1.Cells.Item(103,2).Value = "09:10"
$ws1.Cells.Item(103,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(103,4).Value = 113
$ws1.Cells.Item(103,5).Value = "LP1912"
$ws1.Cells.Item(104,1).Value = "08:31:01"
$ws1.Cells.Item(104,2).Value = "09:11"
$ws1.Cells.Item(104,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(104,4).Value = 40
$ws1.Cells.Item(104,5).Value = "LP1912"
$ws1.Cells.Item(105,1).Value = "07:17:59"
$ws1.Cells.Item(105,2).Value = "09:16"
$ws1.Cells.Item(105,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(105,4).Value = 119
$ws1.Cells.Item(105,5).Value = "LP1912"
$ws1.Cells.Item(106,1).Value = "08:31:01"
$ws1.Cells.Item(106,2).Value = "09:17"
$ws1.Cells.Item(106,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(106,4).Value = 46
$ws1.Cells.Item(106,5).Value = "LP1912"
$ws1.Cells.Item(107,1).Value = "08:31:01"
$ws1.Cells.Item(107,2).Value = "09:19"
$ws1.Cells.Item(107,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(107,4).Value = 48
$ws1.Cells.Item(107,5).Value = "LP1912"
$ws1.Cells.Item(108,1).Value = "08:31:01"
$ws1.Cells.Item(108,2).Value = "09:21"
$ws1.Cells.Item(108,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(108,4).Value = 50
$ws1.Cells.Item(108,5).Value = "LP1912"
$ws1.Cells.Item(109,1).Value = "08:31:01"
$ws1.Cells.Item(109,2).Value = "09:23"
$ws1.Cells.Item(109,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(109,4).Value = 52
$ws1.Cells.Item(109,5).Value = "LP1912"
$ws1.Cells.Item(110,1).Value = "08:31:01"
$ws1.Cells.Item(110,2).Value = "09:23"
$ws1.Cells.Item(110,3).Value = "17_ROMERO"
$ws1.Cells.Item(110,4).Value = 52
$ws1.Cells.Item(110,5).Value = "LP1912"
$ws1.Cells.Item(111,1).Value = "08:31:01"
$ws1.Cells.Item(111,2).Value = "09:24"
$ws1.Cells.Item(111,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(111,4).Value = 53
$ws1.Cells.Item(111,5).Value = "LP1912"

# -- new rows appended at the end --
$ws1.Cells.Item(112,1).Value = "07:48:05"
$ws1.Cells.Item(112,2).Value = "09:32"
$ws1.Cells.Item(112,3).Value = "15_ABASTO"
$ws1.Cells.Item(112,4).Value = 104
$ws1.Cells.Item(112,5).Value = "LP1912"
$ws1.Cells.Item(113,1).Value = "08:31:01"
$ws1.Cells.Item(113,2).Value = "09:33"
$ws1.Cells.Item(113,3).Value = "10_OLMOS"
$ws1.Cells.Item(113,4).Value = 62
$ws1.Cells.Item(113,5).Value = "LP1912"
$ws1.Cells.Item(114,1).Value = "07:48:05"
$ws1.Cells.Item(114,2).Value = "09:34"
$ws1.Cells.Item(114,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(114,4).Value = 106
$ws1.Cells.Item(114,5).Value = "LP1912"
$ws1.Cells.Item(115,1).Value = "08:31:01"
$ws1.Cells.Item(115,2).Value = "09:42"
$ws1.Cells.Item(115,3).Value = "215C_EL PATO"
$ws1.Cells.Item(115,4).Value = 71
$ws1.Cells.Item(115,5).Value = "LP1912"
$ws1.Cells.Item(116,1).Value = "08:31:01"
$ws1.Cells.Item(116,2).Value = "09:44"
$ws1.Cells.Item(116,3).Value = "14_ABASTO"
$ws1.Cells.Item(116,4).Value = 73
$ws1.Cells.Item(116,5).Value = "LP1912"
$ws1.Cells.Item(117,1).Value = "08:31:01"
$ws1.Cells.Item(117,2).Value = "09:52"
$ws1.Cells.Item(117,3).Value = "15_ABASTO"
$ws1.Cells.Item(117,4).Value = 81
$ws1.Cells.Item(117,5).Value = "LP1912"
$ws1.Cells.Item(118,1).Value = "08:31:01"
$ws1.Cells.Item(118,2).Value = "09:56"
$ws1.Cells.Item(118,3).Value = "10_OLMOS"
$ws1.Cells.Item(118,4).Value = 85
$ws1.Cells.Item(118,5).Value = "LP1912"
$ws1.Cells.Item(119,1).Value = "08:31:01"
$ws1.Cells.Item(119,2).Value = "10:11"
$ws1.Cells.Item(119,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(119,4).Value = 100
$ws1.Cells.Item(119,5).Value = "LP1912"
$ws1.Cells.Item(120,1).Value = "08:31:01"
$ws1.Cells.Item(120,2).Value = "10:12"
$ws1.Cells.Item(120,3).Value = "15_ABASTO"
$ws1.Cells.Item(120,4).Value = 101
$ws1.Cells.Item(120,5).Value = "LP1912"
$ws1.Cells.Item(121,1).Value = "08:31:01"
$ws1.Cells.Item(121,2).Value = "10:27"
$ws1.Cells.Item(121,3).Value = "215A_EL PATO"
$ws1.Cells.Item(121,4).Value = 116
$ws1.Cells.Item(121,5).Value = "LP1912"

# ===== Sheet: LP1912-215 =====
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 08:31:01"
$ws2.Range("A3").Value = "Total filas: 19"

# -- changed rows --
$ws2.Cells.Item(22,1).Value = "08:31:01"
$ws2.Cells.Item(22,2).Value = "09:02"
$ws2.Cells.Item(22,3).Value = "215A_EL PATO"
$ws2.Cells.Item(22,4).Value = 31
$ws2.Cells.Item(22,5).Value = "LP1912"
$ws2.Cells.Item(23,1).Value = "08:31:01"
$ws2.Cells.Item(23,2).Value = "09:42"
$ws2.Cells.Item(23,3).Value = "215C_EL PATO"
$ws2.Cells.Item(23,4).Value = 71
$ws2.Cells.Item(23,5).Value = "LP1912"
# -- new row --
$ws2.Cells.Item(24,1).Value = "08:31:01"
$ws2.Cells.Item(24,2).Value = "10:27"
$ws2.Cells.Item(24,3).Value = "215A_EL PATO"
$ws2.Cells.Item(24,4).Value = 116
$ws2.Cells.Item(24,5).Value = "LP1912"

# ===== Sheet: 6203-6173 =====
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 08:31:01"
$ws3.Range("A3").Value = "Total filas: 21"

# -- rows 20-26 rewritten (insert + shift handled by direct overwrite) --
$ws3.Cells.Item(20,1).Value = "08:31:01"
$ws3.Cells.Item(20,2).Value = "08:33"
$ws3.Cells.Item(20,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(20,4).Value = 2
$ws3.Cells.Item(20,5).Value = "L6203"
$ws3.Cells.Item(21,1).Value = "07:17:59"
$ws3.Cells.Item(21,2).Value = "08:34"
$ws3.Cells.Item(21,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(21,4).Value = 77
$ws3.Cells.Item(21,5).Value = "L6173"
$ws3.Cells.Item(22,1).Value = "08:01:10"
$ws3.Cells.Item(22,2).Value = "08:35"
$ws3.Cells.Item(22,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(22,4).Value = 34
$ws3.Cells.Item(22,5).Value = "L6173"
$ws3.Cells.Item(23,1).Value = "08:31:01"
$ws3.Cells.Item(23,2).Value = "08:36"
$ws3.Cells.Item(23,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(23,4).Value = 5
$ws3.Cells.Item(23,5).Value = "L6173"
$ws3.Cells.Item(24,1).Value = "07:17:59"
$ws3.Cells.Item(24,2).Value = "09:08"
$ws3.Cells.Item(24,3).Value = "215D_LA PLATA"
$ws3.Cells.Item(24,4).Value = 111
$ws3.Cells.Item(24,5).Value = "L6203"
$ws3.Cells.Item(25,1).Value = "08:31:01"
$ws3.Cells.Item(25,2).Value = "09:09"
$ws3.Cells.Item(25,3).Value = "215D_LA PLATA"
$ws3.Cells.Item(25,4).Value = 38
$ws3.Cells.Item(25,5).Value = "L6203"
$ws3.Cells.Item(26,1).Value = "08:31:01"
$ws3.Cells.Item(26,2).Value = "10:03"
$ws3.Cells.Item(26,3).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(26,4).Value = 92
$ws3.Cells.Item(26,5).Value = "L6173"
